$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'The applicant has been involved in several projects that align with the job description, such as developing user-friendly platforms using ReactJS and collaborating with back-end developers. The applicant also has experience in developing and fine-tuning models, which showcases their ability to optimize web pages. Although the applicant does not possess NextJS and CSS skills, their proficiency in ReactJS and JavaScript compensates for this gap. The projects demonstrate the applicant''s ability to work with relevant technologies and align with the responsibilities of the Frontend Engineer Intern role.'

$ws.Range("E2").Value = 'The applicant has expressed a strong interest in Japan and the AI field, demonstrating enthusiasm and adaptability. Their willingness to learn a new language, work with a team, and contribute to cutting-edge technologies aligns with the soft skills required by the company. The applicant''s ability to overcome setbacks through collaborative problem-solving and their background in engineering and graphic design showcase adaptability and teamwork. Overall, the applicant exhibits positive attributes and a willingness to embrace new experiences in Japan.'

$ws.Range("C3").Value = 'The applicant has been given a score based on the projects and skills that align well with the job description. The projects ''RAPID'' and ''SMART SENSING MIDDLEWARE'' demonstrate the applicant''s proficiency in frontend development using ReactJS, JavaScript, CSS, and NextJS, which are the required skills for the job. The use of NodeJS, ExpressJS, Socket.IO, WebRTC, and Docker in the ''RAPID'' project also showcases strong backend and real-time web application development skills, indicating the potential for collaborating with backend developers as mentioned in the job description. Although the project ''SHAMIYANA APP'' primarily utilizes Flutter and Dart, the applicant''s overall skills in CSS, Docker, NodeJS, HTML, and JavaScript demonstrate versatility and the ability to adapt to different technologies, reflecting the potential to contribute effectively to frontend engineering. However, the absence of explicit experience with NextJS is a slight gap compared to the required skills for the position, which might have slightly impacted the score.'

$ws.Range("D3").Value = 5

$ws.Range("E3").Value = 'Based on the answers given by the applicant and their demonstrated strengths, the applicant exhibits strong adaptability, teamwork, communication skills, and problem-solving abilities, as evidenced by their willingness to adapt to new cultures, preference for teamwork, and the ability to overcome challenges through collaborative problem-solving. Their clear interest in Japan''s work culture and willingness to learn a new language also indicate a high level of adaptability. Hence, the applicant demonstrates a high level of suitability for the company''s soft skills requirements and willingness to work in Japan.'

$ws.Range("C4").Value = 'The applicant has a score of 75.87 because of their experience in working on projects that involve frontend development using technologies like ReactJS, Django Rest Framework, and Django Rest Framework. The applicant''s projects demonstrate a strong grasp of frontend and backend technologies, along with skills in image processing and machine learning, which are relevant to the job description. However, the applicant lacks direct experience with NextJS, which is a required skill for this position, accounting for the slightly lower score.'

$ws.Range("D4").Value = 5

$ws.Range("E4").Value = 'Based on the applicant''s answers, it is evident that they possess strong communication skills and emphasize a team-oriented approach. Their willingness to adapt to new cultures, demonstrated by their interest in learning a new language in Japan, showcases their adaptability. The applicant''s focus on problem-solving and their ability to learn from setbacks further demonstrates their problem-solving skills. Their willingness to work in Japan and their positive attitude towards embracing a new work culture align with the company''s values, earning them a score of 5.'

$wb.Save()
